$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 440, shifting existing rows 440:537 down to 441:538
$ws.Rows(440).EntireRow.Insert()

# Populate the newly inserted row 440 with the new weekly price record
$ws.Range("A440").Value = 9
$ws.Range("B440").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C440").Value = "Metropolitana"
$ws.Range("D440").Value = 44889
$ws.Range("E440").Value = 13
$ws.Range("F440").Value = "Fruta"
$ws.Range("G440").Value = 100108
$ws.Range("H440").Value = "Tropicales y subtropicales"
$ws.Range("I440").Value = 100108002
$ws.Range("J440").Value = "Mango"
$ws.Range("K440").Value = "Sin especificar"
$ws.Range("L440").Value = "Primera"
$ws.Range("M440").Value = 820
$ws.Range("N440").Value = 7000
$ws.Range("O440").Value = 7500
$ws.Range("P440").Value = 7268
$ws.Range("Q440").Value = "$/bandeja 4 kilos"
$ws.Range("R440").Value = "Brasil"
$ws.Range("S440").Value = 1817
$ws.Range("T440").Value = 4
